$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3:D3").Value = 130
$ws.Range("C5:D5").Value = 32
$ws.Range("C7:D7").Value = 33
$ws.Range("C9:D9").Value = 138
$ws.Range("C11:D11").Value = 105
$ws.Range("C13:D13").Value = 106
$ws.Range("C15:D15").Value = 86
$ws.Range("C17:D17").Value = 36
$ws.Range("C19:D19").Value = 120
$ws.Range("C21:D21").Value = 121
$ws.Range("C23:D23").Value = 112
$ws.Range("C25:D25").Value = 70
$ws.Range("C27:D27").Value = 189
$ws.Range("C29:D29").Value = 131
$ws.Range("C30:D30").Value = 119
$ws.Range("C32:D32").Value = 261
$ws.Range("C34:D34").Value = 88
$ws.Range("C36:D36").Value = 176
$ws.Range("C38:D38").Value = 38
$ws.Range("C40:D40").Value = 117
$ws.Range("C42:D42").Value = 41
$ws.Range("C44:D44").Value = 92
$ws.Range("C46:D46").Value = 95
$ws.Range("C48:D48").Value = 149
$ws.Range("C50:D50").Value = 73
$ws.Range("C52:D52").Value = 69
$ws.Range("C54:D54").Value = 155
$ws.Range("C56:D56").Value = 693
$ws.Range("C58:D58").Value = 118
$ws.Range("C60:D60").Value = 76
$ws.Range("C62:D62").Value = 83
$ws.Range("C64:D64").Value = 72
$ws.Range("C66:D66").Value = 78
$ws.Range("C68:D68").Value = 124
$ws.Range("C70:D70").Value = 114
$ws.Range("C72:D72").Value = 62
$ws.Range("C74:D74").Value = 45
$ws.Range("C76:D76").Value = 80
$ws.Range("C77").Value = 116.5
